$wb = $excel.ActiveWorkbook

# The second sheet is currently named "9358" and needs to be renamed to "8059"
$ws = $wb.Worksheets.Item(2)
$ws.Name = "8059"

# Update unit number cell to match the new sheet name (force text so it
# doesn't get interpreted as a number)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "8059"

# SAMPLES DUE (Y/N): set to Y
$ws.Range("F3").Value = "Y"

# FRA DEAD DATE changes (force text so the date-like string isn't
# converted into a date serial number)
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "01-06-2020"

# LUBE DUE (Y/N): set to Y
$ws.Range("F5").Value = "Y"

# MI DUE (month): set to 6mo
$ws.Range("C6").Value = "6mo"

# C/S DUE (Y/N): set to Y
$ws.Range("F6").Value = "Y"

# AF DUE (Y/N): set to Y
$ws.Range("F7").Value = "Y"

# Notes section updates
$ws.Range("B24").Value = "AIR BRAKE"
$ws.Range("B25").Value = "Alertor penalty, source still present"
